$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) column C for all existing data rows (2-509)
#    from 2023-09-21 (45190) to 2023-09-23 (45192).
$ws.Range("C2:C509").Value = 45192

# 2. Add new row 510: A 44463-2023
$ws.Range("A510").Value = "A 44463-2023"
$ws.Range("B510").NumberFormat = "YYYY-MM-DD"
$ws.Range("B510").Value = 45189
$ws.Range("C510").NumberFormat = "YYYY-MM-DD"
$ws.Range("C510").Value = 45192
$ws.Range("D510").Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E510").Value = "ULRICEHAMN"
$ws.Range("G510").Value = 2.1
$ws.Range("H510").Value = 0
$ws.Range("I510").Value = 0
$ws.Range("J510").Value = 0
$ws.Range("K510").Value = 0
$ws.Range("L510").Value = 0
$ws.Range("M510").Value = 0
$ws.Range("N510").Value = 0
$ws.Range("O510").Value = 0
$ws.Range("P510").Value = 0
$ws.Range("Q510").Value = 0
$ws.Range("R510").WrapText = $true
$ws.Range("R510").Value = ""
$ws.Range("A510").RowHeight = 15

# 3. Add new row 511: A 44930-2023
$ws.Range("A511").Value = "A 44930-2023"
$ws.Range("B511").NumberFormat = "YYYY-MM-DD"
$ws.Range("B511").Value = 45190
$ws.Range("C511").NumberFormat = "YYYY-MM-DD"
$ws.Range("C511").Value = 45192
$ws.Range("D511").Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E511").Value = "ULRICEHAMN"
$ws.Range("G511").Value = 5.1
$ws.Range("H511").Value = 0
$ws.Range("I511").Value = 0
$ws.Range("J511").Value = 0
$ws.Range("K511").Value = 0
$ws.Range("L511").Value = 0
$ws.Range("M511").Value = 0
$ws.Range("N511").Value = 0
$ws.Range("O511").Value = 0
$ws.Range("P511").Value = 0
$ws.Range("Q511").Value = 0
$ws.Range("R511").WrapText = $true
$ws.Range("R511").Value = ""

# 4. Ensure row 509's height is explicitly set to 15 (customHeight), matching the
#    author's re-save of that row.
$ws.Range("A509").RowHeight = 15
